# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the Leve profit tables (columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 325.66666
$ws.Range("I2").Value = 245.8
$ws.Range("J2").Value = 725
$ws.Range("K2").Value = 245.8
$ws.Range("L2").Value = 725
$ws.Range("M2").Value = -132.8
$ws.Range("N2").Value = -951

$ws.Range("H51").Value = 3908.3333
$ws.Range("J51").Value = 3990
$ws.Range("L51").Value = 3990
$ws.Range("N51").Value = -4958

$ws.Range("H64").Value = 4938
$ws.Range("I64").Value = 4876
$ws.Range("K64").Value = 4876
$ws.Range("M64").Value = -4628

$ws.Range("H67").Value = 4938
$ws.Range("I67").Value = 4876
$ws.Range("K67").Value = 4876
$ws.Range("M67").Value = -4018

$ws.Range("H116").Value = 987169.7
$ws.Range("I116").Value = 6533.5
$ws.Range("K116").Value = 6533.5
$ws.Range("M116").Value = -3091.5

$ws.Range("H137").Value = 519713.53
$ws.Range("I137").Value = 1568.75
$ws.Range("K137").Value = 4706.25
$ws.Range("M137").Value = -2156.25

$ws.Range("H138").Value = 2235.1667
$ws.Range("I138").Value = 946.4211
$ws.Range("K138").Value = 2839.2633
$ws.Range("M138").Value = 2300.7367

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8567.359
$ws.Range("I32").Value = 3473.6042
$ws.Range("J32").Value = 23848.625
$ws.Range("K32").Value = 3473.6042
$ws.Range("L32").Value = 23848.625
$ws.Range("M32").Value = -3186.6042
$ws.Range("N32").Value = -24422.625

$ws.Range("H61").Value = 1811.2727
$ws.Range("I61").Value = 1303
$ws.Range("J61").Value = 3166.6667
$ws.Range("K61").Value = 1303
$ws.Range("L61").Value = 3166.6667
$ws.Range("M61").Value = -1091
$ws.Range("N61").Value = -3590.6667

$ws.Range("H122").Value = 4455.2085
$ws.Range("I122").Value = 4842.9443
$ws.Range("J122").Value = 3292
$ws.Range("K122").Value = 14528.8329
$ws.Range("L122").Value = 9876
$ws.Range("M122").Value = -12078.8329
$ws.Range("N122").Value = -14776

$ws.Range("H132").Value = 2357.5417
$ws.Range("I132").Value = 1806.0667
$ws.Range("J132").Value = 3276.6667
$ws.Range("K132").Value = 5418.2001
$ws.Range("L132").Value = 9830.000100000001
$ws.Range("M132").Value = -2888.2001
$ws.Range("N132").Value = -14890.0001

$ws.Range("H136").Value = 1811.2727
$ws.Range("I136").Value = 1303
$ws.Range("J136").Value = 3166.6667
$ws.Range("K136").Value = 3909
$ws.Range("L136").Value = 9500.000100000001
$ws.Range("M136").Value = -1359
$ws.Range("N136").Value = -14600.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 44993.332
$ws.Range("J13").Value = 44993.332
$ws.Range("L13").Value = 44993.332
$ws.Range("N13").Value = -45329.332

$ws.Range("H22").Value = 87173.086
$ws.Range("I22").Value = 111786.336
$ws.Range("K22").Value = 111786.336
$ws.Range("M22").Value = -111613.336

$ws.Range("H132").Value = 27188.781
$ws.Range("J132").Value = 27188.781
$ws.Range("L132").Value = 27188.781
$ws.Range("N132").Value = -37308.781

$ws.Range("H134").Value = 1854.4814
$ws.Range("I134").Value = 1274.9524
$ws.Range("K134").Value = 3824.857199999999
$ws.Range("M134").Value = -1289.857199999999

$ws.Range("H140").Value = 90561.42999999999
$ws.Range("J140").Value = 99988.336
$ws.Range("L140").Value = 99988.336
$ws.Range("N140").Value = -110348.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H58").Value = 1610.5
$ws.Range("I58").Value = 1524.9286
$ws.Range("K58").Value = 1524.9286
$ws.Range("M58").Value = -1321.9286

$ws.Range("H59").Value = 140000
$ws.Range("J59").Value = 140000
$ws.Range("L59").Value = 140000
$ws.Range("N59").Value = -142290

$ws.Range("H107").Value = 811
$ws.Range("I107").Value = 827.6923
$ws.Range("J107").Value = 702.5
$ws.Range("K107").Value = 827.6923
$ws.Range("L107").Value = 702.5
$ws.Range("M107").Value = 1092.3077
$ws.Range("N107").Value = -4542.5

$ws.Range("H132").Value = 2672.125
$ws.Range("I132").Value = 2479.6667
$ws.Range("K132").Value = 7439.000100000001
$ws.Range("M132").Value = -4909.000100000001

$ws.Range("H134").Value = 2527.8333
$ws.Range("I134").Value = 2217.3635
$ws.Range("J134").Value = 3015.7144
$ws.Range("K134").Value = 6652.0905
$ws.Range("L134").Value = 9047.143199999999
$ws.Range("M134").Value = -4117.0905
$ws.Range("N134").Value = -14117.1432

$ws.Range("H136").Value = 1610.5
$ws.Range("I136").Value = 1524.9286
$ws.Range("K136").Value = 4574.7858
$ws.Range("M136").Value = -2024.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 94.117645
$ws.Range("J2").Value = 100.44444
$ws.Range("L2").Value = 602.66664
$ws.Range("N2").Value = -828.66664

$ws.Range("H57").Value = 2498.75
$ws.Range("I57").Value = 500
$ws.Range("K57").Value = 1500
$ws.Range("M57").Value = -941

$ws.Range("H59").Value = 1398.8889
$ws.Range("I59").Value = 120
$ws.Range("J59").Value = 2997.5
$ws.Range("K59").Value = 360
$ws.Range("L59").Value = 8992.5
$ws.Range("M59").Value = 180
$ws.Range("N59").Value = -10072.5

$ws.Range("H81").Value = 1484
$ws.Range("I81").Value = 777.5
$ws.Range("J81").Value = 2897
$ws.Range("K81").Value = 2332.5
$ws.Range("L81").Value = 8691
$ws.Range("M81").Value = -1209.5
$ws.Range("N81").Value = -10937

$ws.Range("H84").Value = 1484
$ws.Range("I84").Value = 777.5
$ws.Range("J84").Value = 2897
$ws.Range("K84").Value = 6997.5
$ws.Range("L84").Value = 26073
$ws.Range("M84").Value = -1381.5
$ws.Range("N84").Value = -37305

$ws.Range("H139").Value = 7707.9614
$ws.Range("I139").Value = 4019.5
$ws.Range("J139").Value = 10869.5
$ws.Range("K139").Value = 12058.5
$ws.Range("L139").Value = 32608.5
$ws.Range("M139").Value = -6918.5
$ws.Range("N139").Value = -42888.5

$ws.Range("H140").Value = 2063.4
$ws.Range("I140").Value = 1239
$ws.Range("K140").Value = 3717
$ws.Range("M140").Value = 1463

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 300
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -746

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 46134.09
$ws.Range("J7").Value = 56994.8
$ws.Range("L7").Value = 56994.8
$ws.Range("N7").Value = -57218.8

$ws.Range("H22").Value = 1150.5
$ws.Range("J22").Value = 1204.1666
$ws.Range("L22").Value = 1204.1666
$ws.Range("N22").Value = -1794.1666

$ws.Range("H27").Value = 1150.5
$ws.Range("J27").Value = 1204.1666
$ws.Range("L27").Value = 1204.1666
$ws.Range("N27").Value = -1418.1666

$ws.Range("H46").Value = 3279.6667
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3279.6667
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3655.6667

$ws.Range("H68").Value = 1051750.5
$ws.Range("I68").Value = 1051750.5
$ws.Range("K68").Value = 1051750.5
$ws.Range("M68").Value = -1051001.5

$ws.Range("H71").Value = 1051750.5
$ws.Range("I71").Value = 1051750.5
$ws.Range("K71").Value = 5258752.5
$ws.Range("M71").Value = -5255008.5

$ws.Range("H100").Value = 29375.75
$ws.Range("I100").Value = 29375.75
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 29375.75
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H126").Value = 46134.09
$ws.Range("J126").Value = 56994.8
$ws.Range("L126").Value = 170984.4
$ws.Range("N126").Value = -175924.4

$ws.Range("H132").Value = 10593.667
$ws.Range("I132").Value = 15870.6
$ws.Range("J132").Value = 3997.5
$ws.Range("K132").Value = 47611.8
$ws.Range("L132").Value = 11992.5
$ws.Range("M132").Value = -45081.8
$ws.Range("N132").Value = -17052.5

$ws.Range("H140").Value = 71143
$ws.Range("J140").Value = 74214.5
$ws.Range("L140").Value = 74214.5
$ws.Range("N140").Value = -84574.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 20000
$ws.Range("K51").Value = 20000
$ws.Range("M51").Value = -19490

$ws.Range("H52").Value = 32500
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 32500
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -32952

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H122").Value = 2660
$ws.Range("I122").Value = 1706.2307
$ws.Range("J122").Value = 5139.8
$ws.Range("K122").Value = 5118.6921
$ws.Range("L122").Value = 15419.4
$ws.Range("M122").Value = -2668.6921
$ws.Range("N122").Value = -20319.4
